$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Replace the embedded "WCP.jpg" picture with a hyperlink that
#    points at the image's original URL (the picture is removed and
#    a plain-text hyperlink run is inserted in its place, inside the
#    same BodyText paragraph that used to hold the drawing).
# ------------------------------------------------------------------
$picUrl = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/WCP.jpg"

if ($d.InlineShapes.Count -ge 1) {
    $shape = $d.InlineShapes(1)
    $shapeRange = $shape.Range
    $shape.Delete()
    $null = $d.Hyperlinks.Add($shapeRange, $picUrl, $null, $null, $picUrl)
}

Write-Output "done"
